# Generate Report for handoff
# Update "Latest Handoff Datetime" for the file that is ready for handoff
# (6172aa8c-fd0c-45bd-b147-5e3b6aba0a9e) on each locale sheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    "zh-cn" = "2016-01-13 04:08:25"
    "de-de" = "2016-01-13 04:08:46"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newDatetime = $updates[$sheetName]

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 1; $r -le $lastRow; $r++) {
        $sourceName = $ws.Cells.Item($r, 1).Value2
        if ($sourceName -like "6172aa8c-fd0c-45bd-b147-5e3b6aba0a9e*") {
            $ws.Cells.Item($r, 4).Value = $newDatetime
            break
        }
    }
}
